# Apply data-cleanup edit:
#  - rename header "player_id" -> "player_id_x"
#  - rename header "birth_year" -> "birth_year_x"
#  - fix player_id values in column C (rows 2-23) from 3463 -> 3462

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1)
$ws.Range("C1").Value = "player_id_x"
$ws.Range("E1").Value = "birth_year_x"

# Correct player_id values for rows 2 through 23 (Michael Jordan rows)
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 3).Value = 3462
}
